$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.864.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.532.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.40'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.43'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.530.42'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.96%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.485'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.47%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.57%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.12'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.81%  '

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.124.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.545.81'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.891.58'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.61%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.49'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '451.55'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.639'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.93'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.672.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.90%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.23'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.35'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -7.62%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.69'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.29%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.55'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.95'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.90'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.20'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.158'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.524.09'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.11'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.93'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.20%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.61'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.89%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0876'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.892'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.79'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.36'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.73%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.38%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.66'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.48%  '

